# Generate Report for Archive
#
# 1. Replace the "Ready for handoff" status text with "In Translation"
#    everywhere it appears (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 -
#    these all share the same shared-string entry in the source file).
# 2. Shrink the now-narrower status columns to match the new text:
#    Overview columns E & F, and column C on the zh-cn / de-de sheets.

$wb  = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the status text -------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- Resize the status columns to fit the shorter text ----------------------
# COM only accepts ColumnWidth in character units (quantized to whole
# pixels internally), so pick the value whose stored width lands on the
# new target width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
